# "Updated symbol list" run — refresh the coin price/volume snapshot on the
# "cryptos" sheet. Prices (column D) are stored as literal text in this sheet
# (the site's scraper writes them as plain strings, e.g. "248.45", not
# numbers), so every numeric-looking update below is entered with a leading
# apostrophe to force text entry and keep the exact digits/trailing zeros
# instead of letting Excel coerce it to a floating point number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- simple price (column D) refreshes for the top of the table ---
$ws.Range("D2").Value = '''248.39'
$ws.Range("D3").Value = '''21.72'
$ws.Range("D4").Value = '''5.443'
$ws.Range("D5").Value = '''0.05688'
$ws.Range("D6").Value = '''3.384'
$ws.Range("D7").Value = '''0.8071'
$ws.Range("D8").Value = '''1.028'

# --- rows 9-17: the coin ranked #9 (One/ONE) dropped out of the top ranks
# and off the bottom of this block, so every coin shifted up one row
# (row 10's coin -> row 9, row 11's coin -> row 10, ...), with "One" moving
# from row 9 down to row 17 (now flagged "Worstin24h" instead of
# "Bestin24h"). Each row gets the next row's old name/link plus a freshly
# scraped price, and the rank-prefixed "label" in column E. ---
$ws.Range("B9").Value = 'WazirX'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D9").Value = '''0.1458'
$ws.Range("E9").Value = '8WazirXWRX'

$ws.Range("B10").Value = 'MandalaExchangeToken'
$ws.Range("C10").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D10").Value = '''0.07634'
$ws.Range("E10").Value = '9MandalaExchangeTokenMDX'

$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D11").Value = '''0.03165'
$ws.Range("E11").Value = '10LiechtensteinCryptoassetsExchangeLCX'

$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = '''0.03029'
$ws.Range("E12").Value = '11BitrueCoinBTR'

$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = '''0.09261'
$ws.Range("E13").Value = '12BitMartTokenBMX'

$ws.Range("B14").Value = 'MCDex'
$ws.Range("C14").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D14").Value = '''3.520'
$ws.Range("E14").Value = '13MCDexMCB'

$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = '''0.001650'
$ws.Range("E15").Value = '14BitForexTokenBF'

$ws.Range("B16").Value = 'CoinExToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D16").Value = '''0.04708'
$ws.Range("E16").Value = '15CoinExTokenCET'

$ws.Range("B17").Value = 'One'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D17").Value = '''0.0005863'
$ws.Range("E17").Value = '16OneONEWorstin24h'

# --- remaining scattered price (and a couple of "Bestin24h"/"Worstin24h"
# label) refreshes further down the table ---
$ws.Range("D18").Value = '''0.006347'

$ws.Range("D19").Value = '''0.005028'
$ws.Range("E19").Value = '18HotbitTokenHTBBestin24h'

$ws.Range("D21").Value = '''0.0001501'
$ws.Range("D22").Value = '''0.0003102'
$ws.Range("D24").Value = '''6.417'
$ws.Range("D25").Value = '''2.175'
$ws.Range("D26").Value = '''0.3296'
$ws.Range("D40").Value = '''0.04069'
$ws.Range("D41").Value = '''0.006984'
$ws.Range("D42").Value = '''0.003503'
$ws.Range("D43").Value = '''0.1042'
$ws.Range("D44").Value = '''0.007878'
$ws.Range("D45").Value = '''0.00005908'

$ws.Range("D47").Value = '''0.0005503'
$ws.Range("E47").Value = '46ACDXExchangeACXT'

$ws.Range("D48").Value = '''0.6828'
$ws.Range("D49").Value = '''0.008045'
